$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142, shifting existing rows 142-192 down to 143-193
$ws.Rows.Item(142).EntireRow.Insert()

# Populate the newly inserted row 142 with the new data record
$ws.Cells.Item(142, 1).Value = 3
$ws.Cells.Item(142, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(142, 3).Value = "Coquimbo"
$ws.Cells.Item(142, 4).Value = 44809
$ws.Cells.Item(142, 5).Value = 5
$ws.Cells.Item(142, 6).Value = 100112026
$ws.Cells.Item(142, 7).Value = "Haba"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 85
$ws.Cells.Item(142, 11).Value = 13000
$ws.Cells.Item(142, 12).Value = 14000
$ws.Cells.Item(142, 13).Value = 13471
$ws.Cells.Item(142, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(142, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(142, 16).Value = 539
$ws.Cells.Item(142, 17).Value = 25
$ws.Cells.Item(142, 18).Value = "Hortaliza"
